$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, preserving its existing (unstyled) cell format
# so that numeric-looking strings (e.g. "0.995") are not reinterpreted as numbers.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "27.783.71"
Set-TextValue $ws.Range("E2") "  -0.57%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.627.81"
Set-TextValue $ws.Range("E3") "  -1.16%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.995"
Set-TextValue $ws.Range("E4") "  -0.59%  "

# Row 5
Set-TextValue $ws.Range("D5") "211.12"
Set-TextValue $ws.Range("E5") "  -1.18%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -1.10%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.994"
Set-TextValue $ws.Range("E7") "  -0.61%  "

# Row 8
Set-TextValue $ws.Range("D8") "23.24"
Set-TextValue $ws.Range("E8") "  -0.98%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.257"
Set-TextValue $ws.Range("E9") "  -2.94%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0613"
Set-TextValue $ws.Range("E10") "  -0.41%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0880"
Set-TextValue $ws.Range("E11") "  +0.87%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.859.41"
Set-TextValue $ws.Range("E12") "  -1.13%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.635.13"
Set-TextValue $ws.Range("E13") "  -0.79%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.04"
Set-TextValue $ws.Range("E14") "  -0.89%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.563"
Set-TextValue $ws.Range("E15") "  -0.43%  "

# Row 16
Set-TextValue $ws.Range("D16") "65.15"
Set-TextValue $ws.Range("E16") "  -0.76%  "

# Row 17
Set-TextValue $ws.Range("D17") "27.801.90"
Set-TextValue $ws.Range("E17") "  -0.60%  "

# Row 18
Set-TextValue $ws.Range("D18") "230.00"
Set-TextValue $ws.Range("E18") "  -0.81%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.0₃0721"
Set-TextValue $ws.Range("E19") "  -0.32%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.51"
Set-TextValue $ws.Range("E20") "  -2.04%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.995"
Set-TextValue $ws.Range("E21") "  -0.61%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.35"
Set-TextValue $ws.Range("E22") "  -0.95%  "

# Row 23
Set-TextValue $ws.Range("D23") "10.30"
Set-TextValue $ws.Range("E23") "  -3.49%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.05"
Set-TextValue $ws.Range("E24") "  -4.22%  "

# Row 25
Set-TextValue $ws.Range("D25") "154.02"
Set-TextValue $ws.Range("E25") "  +1.12%  "

# Row 26
Set-TextValue $ws.Range("D26") "6.94"
Set-TextValue $ws.Range("E26") "  +0.36%  "

# Row 27
Set-TextValue $ws.Range("B27") "Stellar"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D27") "0.111"
Set-TextValue $ws.Range("E27") "  -1.23%  "

# Row 28
Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "15.60"
Set-TextValue $ws.Range("E28") "  -1.02%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.996"
Set-TextValue $ws.Range("E29") "  -0.55%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -1.25%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.0480"
Set-TextValue $ws.Range("E31") "  -1.01%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +1.74%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.401.95"
Set-TextValue $ws.Range("E33") "  -2.80%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.07"
Set-TextValue $ws.Range("E34") "  -0.09%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -0.16%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.01"
Set-TextValue $ws.Range("E36") "  +8.30%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.35"
Set-TextValue $ws.Range("E37") "  +0.76%  "

# Row 38
Set-TextValue $ws.Range("E38") "  +0.13%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.560"
Set-TextValue $ws.Range("E39") "  +0.23%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.867"
Set-TextValue $ws.Range("E40") "  -2.44%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -0.31%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.995"
Set-TextValue $ws.Range("E42") "  -0.57%  "

# Row 43
Set-TextValue $ws.Range("D43") "66.57"
Set-TextValue $ws.Range("E43") "  -3.73%  "

# Row 44
Set-TextValue $ws.Range("D44") "5.50"
Set-TextValue $ws.Range("E44") "  +1.43%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.82"
Set-TextValue $ws.Range("E45") "  -0.31%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.19"
Set-TextValue $ws.Range("E46") "  -1.43%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.768.93"
Set-TextValue $ws.Range("E47") "  -1.18%  "

# Row 48
Set-TextValue $ws.Range("D48") "87.65"
Set-TextValue $ws.Range("E48") "  -1.67%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0996"
Set-TextValue $ws.Range("E49") "  -1.47%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0505"
Set-TextValue $ws.Range("E50") "  -0.57%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0₇0977"
Set-TextValue $ws.Range("E51") "  -6.46%  "

